$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H17").Value = 2948768
$ws.Range("J17").Value = 2948768
$ws.Range("L17").Value = 8846304
$ws.Range("N17").Value = -8846640
$ws.Range("H40").Value = 1811
$ws.Range("I40").Value = 1685.6666
$ws.Range("K40").Value = 1685.6666
$ws.Range("M40").Value = -1510.6666
$ws.Range("H70").Value = 815.6
$ws.Range("I70").Value = 802.72
$ws.Range("J70").Value = 880
$ws.Range("K70").Value = 2408.16
$ws.Range("L70").Value = 2640
$ws.Range("M70").Value = -2138.16
$ws.Range("N70").Value = -3180
$ws.Range("H73").Value = 815.6
$ws.Range("I73").Value = 802.72
$ws.Range("J73").Value = 880
$ws.Range("K73").Value = 2408.16
$ws.Range("L73").Value = 2640
$ws.Range("M73").Value = -1472.16
$ws.Range("N73").Value = -4512
$ws.Range("H76").Value = 3385.6
$ws.Range("I76").Value = 2831.25
$ws.Range("J76").Value = 5603
$ws.Range("K76").Value = 2831.25
$ws.Range("L76").Value = 5603
$ws.Range("M76").Value = -2516.25
$ws.Range("N76").Value = -6233
$ws.Range("H79").Value = 3385.6
$ws.Range("I79").Value = 2831.25
$ws.Range("J79").Value = 5603
$ws.Range("K79").Value = 2831.25
$ws.Range("L79").Value = 5603
$ws.Range("M79").Value = -1739.25
$ws.Range("N79").Value = -7787
$ws.Range("H101").Value = 729.5
$ws.Range("I101").Value = 729.5
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 2188.5
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = -566.5
$ws.Range("N101").ClearContents()
$ws.Range("H126").Value = 32792.5
$ws.Range("J126").Value = 32792.5
$ws.Range("L126").Value = 32792.5
$ws.Range("N126").Value = -42672.5
$ws.Range("H129").Value = 583.6875
$ws.Range("I129").Value = 479.92307
$ws.Range("J129").Value = 1033.3334
$ws.Range("K129").Value = 1439.76921
$ws.Range("L129").Value = 3100.0002
$ws.Range("M129").Value = 3560.23079
$ws.Range("N129").Value = -13100.0002
$ws.Range("H130").Value = 15404.737
$ws.Range("J130").Value = 15999.444
$ws.Range("L130").Value = 15999.444
$ws.Range("N130").Value = -26039.444
$ws.Range("H138").Value = 7383.0156
$ws.Range("I138").Value = 3980.7273
$ws.Range("J138").Value = 8089.151
$ws.Range("K138").Value = 11942.1819
$ws.Range("L138").Value = 24267.453
$ws.Range("M138").Value = -6802.1819
$ws.Range("N138").Value = -34547.453

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2186.0454
$ws.Range("I45").Value = 1569.625
$ws.Range("K45").Value = 1569.625
$ws.Range("M45").Value = -1192.625
$ws.Range("H74").Value = 1567.873
$ws.Range("I74").Value = 1093.5
$ws.Range("J74").Value = 3810.3635
$ws.Range("K74").Value = 1093.5
$ws.Range("L74").Value = 3810.3635
$ws.Range("M74").Value = -219.5
$ws.Range("N74").Value = -5558.363499999999
$ws.Range("H77").Value = 1567.873
$ws.Range("I77").Value = 1093.5
$ws.Range("J77").Value = 3810.3635
$ws.Range("K77").Value = 5467.5
$ws.Range("L77").Value = 19051.8175
$ws.Range("M77").Value = -1099.5
$ws.Range("N77").Value = -27787.8175
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H122").Value = 2838
$ws.Range("I122").Value = 2250
$ws.Range("J122").Value = 4014
$ws.Range("K122").Value = 6750
$ws.Range("L122").Value = 12042
$ws.Range("M122").Value = -4300
$ws.Range("N122").Value = -16942
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 4173.3335
$ws.Range("I7").Value = 260
$ws.Range("J7").Value = 12000
$ws.Range("K7").Value = 260
$ws.Range("L7").Value = 12000
$ws.Range("M7").Value = -147
$ws.Range("N7").Value = -12226
$ws.Range("H94").Value = 2155.3635
$ws.Range("I94").Value = 2297.8
$ws.Range("J94").Value = 2036.6666
$ws.Range("K94").Value = 2297.8
$ws.Range("L94").Value = 2036.6666
$ws.Range("M94").Value = -1846.8
$ws.Range("N94").Value = -2938.6666
$ws.Range("H134").Value = 2044.625
$ws.Range("I134").Value = 1610.3636
$ws.Range("K134").Value = 4831.0908
$ws.Range("M134").Value = -2296.0908

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9729.619000000001
$ws.Range("I31").Value = 3808.8572
$ws.Range("J31").Value = 21571.143
$ws.Range("K31").Value = 3808.8572
$ws.Range("L31").Value = 21571.143
$ws.Range("M31").Value = -3513.8572
$ws.Range("N31").Value = -22161.143
$ws.Range("H34").Value = 9729.619000000001
$ws.Range("I34").Value = 3808.8572
$ws.Range("J34").Value = 21571.143
$ws.Range("K34").Value = 3808.8572
$ws.Range("L34").Value = 21571.143
$ws.Range("M34").Value = -3606.8572
$ws.Range("N34").Value = -21975.143
$ws.Range("H132").Value = 4328.4287
$ws.Range("I132").Value = 1900
$ws.Range("J132").Value = 5299.8
$ws.Range("K132").Value = 5700
$ws.Range("L132").Value = 15899.4
$ws.Range("M132").Value = -3170
$ws.Range("N132").Value = -20959.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 3999
$ws.Range("J35").Value = 3999
$ws.Range("L35").Value = 11997
$ws.Range("N35").Value = -12573
$ws.Range("H122").Value = 1063.3214
$ws.Range("I122").Value = 611.2083
$ws.Range("J122").Value = 3776
$ws.Range("K122").Value = 5500.8747
$ws.Range("L122").Value = 33984
$ws.Range("M122").Value = -3050.8747
$ws.Range("N122").Value = -38884
$ws.Range("H131").Value = 1056.04
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 1056.04
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 3168.12
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -13248.12
$ws.Range("H138").Value = 1815.3334
$ws.Range("I138").Value = 1723
$ws.Range("J138").Value = 2000
$ws.Range("K138").Value = 5169
$ws.Range("L138").Value = 6000
$ws.Range("M138").Value = -29
$ws.Range("N138").Value = -16280

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 12750
$ws.Range("J5").Value = 14600
$ws.Range("L5").Value = 14600
$ws.Range("N5").Value = -14824
$ws.Range("H122").Value = 2418
$ws.Range("I122").Value = 1879.7084
$ws.Range("J122").Value = 3340.7856
$ws.Range("K122").Value = 5639.1252
$ws.Range("L122").Value = 10022.3568
$ws.Range("M122").Value = -3189.1252
$ws.Range("N122").Value = -14922.3568

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 74657.78999999999
$ws.Range("I7").Value = 102940.4
$ws.Range("K7").Value = 102940.4
$ws.Range("M7").Value = -102828.4
$ws.Range("H40").Value = 87000
$ws.Range("I40").Value = 103600
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 103600
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -103464
$ws.Range("N40").Value = -4272
$ws.Range("H122").Value = 18523638
$ws.Range("J122").Value = 2680
$ws.Range("L122").Value = 8040
$ws.Range("N122").Value = -12940
$ws.Range("H126").Value = 74657.78999999999
$ws.Range("I126").Value = 102940.4
$ws.Range("K126").Value = 308821.2
$ws.Range("M126").Value = -306351.2

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 621.41174
$ws.Range("I107").Value = 614.1539
$ws.Range("J107").Value = 645
$ws.Range("K107").Value = 1842.4617
$ws.Range("L107").Value = 1935
$ws.Range("M107").Value = 77.53829999999994
$ws.Range("N107").Value = -5775
$ws.Range("H122").Value = 39537.85
$ws.Range("I122").Value = 64707.688
$ws.Range("J122").Value = 2927.182
$ws.Range("K122").Value = 194123.064
$ws.Range("L122").Value = 8781.545999999998
$ws.Range("M122").Value = -191673.064
$ws.Range("N122").Value = -13681.546
$ws.Range("H132").Value = 2605.2327
$ws.Range("I132").Value = 1873.9231
$ws.Range("J132").Value = 3723.7058
$ws.Range("K132").Value = 5621.7693
$ws.Range("L132").Value = 11171.1174
$ws.Range("M132").Value = -3091.7693
$ws.Range("N132").Value = -16231.1174
